{"js": "// Replace the 25 two-digit-division answer strings with their updated\n// values, as described by the diff (content-only text swap inside\n// existing table cells; formatting/runs are untouched).\nconst replacements = [\n  [\"27\u00f78=3, 3\", \"66\u00f72=33, 0\"],\n  [\"78\u00f75=15, 3\", \"97\u00f76=16, 1\"],\n  [\"55\u00f75=11, 0\", \"57\u00f74=14, 1\"],\n  [\"84\u00f72=42, 0\", \"91\u00f73=30, 1\"],\n  [\"47\u00f75=9, 2\", \"31\u00f74=7, 3\"],\n  [\"51\u00f78=6, 3\", \"72\u00f76=12, 0\"],\n  [\"67\u00f79=7, 4\", \"87\u00f74=21, 3\"],\n  [\"28\u00f77=4, 0\", \"30\u00f73=10, 0\"],\n  [\"24\u00f72=12, 0\", \"21\u00f72=10, 1\"],\n  [\"62\u00f72=31, 0\", \"28\u00f78=3, 4\"],\n  [\"22\u00f76=3, 4\", \"57\u00f75=11, 2\"],\n  [\"50\u00f77=7, 1\", \"24\u00f73=8, 0\"],\n  [\"77\u00f77=11, 0\", \"95\u00f78=11, 7\"],\n  [\"41\u00f76=6, 5\", \"92\u00f73=30, 2\"],\n  [\"13\u00f79=1, 4\", \"76\u00f72=38, 0\"],\n  [\"43\u00f76=7, 1\", \"93\u00f78=11, 5\"],\n  [\"32\u00f74=8, 0\", \"59\u00f77=8, 3\"],\n  [\"36\u00f79=4, 0\", \"61\u00f78=7, 5\"],\n  [\"24\u00f74=6, 0\", \"34\u00f78=4, 2\"],\n  [\"83\u00f78=10, 3\", \"82\u00f76=13, 4\"],\n  [\"95\u00f77=13, 4\", \"14\u00f72=7, 0\"],\n  [\"19\u00f78=2, 3\", \"38\u00f77=5, 3\"],\n  [\"34\u00f79=3, 7\", \"71\u00f77=10, 1\"],\n  [\"46\u00f74=11, 2\", \"83\u00f76=13, 5\"],\n  [\"68\u00f74=17, 0\", \"16\u00f78=2, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-division answer strings with their updated\n# values, as described by the diff (content-only text swap inside\n# existing table cells; formatting/runs are untouched).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"27\u00f78=3, 3\", \"66\u00f72=33, 0\"),\n    @(\"78\u00f75=15, 3\", \"97\u00f76=16, 1\"),\n    @(\"55\u00f75=11, 0\", \"57\u00f74=14, 1\"),\n    @(\"84\u00f72=42, 0\", \"91\u00f73=30, 1\"),\n    @(\"47\u00f75=9, 2\", \"31\u00f74=7, 3\"),\n    @(\"51\u00f78=6, 3\", \"72\u00f76=12, 0\"),\n    @(\"67\u00f79=7, 4\", \"87\u00f74=21, 3\"),\n    @(\"28\u00f77=4, 0\", \"30\u00f73=10, 0\"),\n    @(\"24\u00f72=12, 0\", \"21\u00f72=10, 1\"),\n    @(\"62\u00f72=31, 0\", \"28\u00f78=3, 4\"),\n    @(\"22\u00f76=3, 4\", \"57\u00f75=11, 2\"),\n    @(\"50\u00f77=7, 1\", \"24\u00f73=8, 0\"),\n    @(\"77\u00f77=11, 0\", \"95\u00f78=11, 7\"),\n    @(\"41\u00f76=6, 5\", \"92\u00f73=30, 2\"),\n    @(\"13\u00f79=1, 4\", \"76\u00f72=38, 0\"),\n    @(\"43\u00f76=7, 1\", \"93\u00f78=11, 5\"),\n    @(\"32\u00f74=8, 0\", \"59\u00f77=8, 3\"),\n    @(\"36\u00f79=4, 0\", \"61\u00f78=7, 5\"),\n    @(\"24\u00f74=6, 0\", \"34\u00f78=4, 2\"),\n    @(\"83\u00f78=10, 3\", \"82\u00f76=13, 4\"),\n    @(\"95\u00f77=13, 4\", \"14\u00f72=7, 0\"),\n    @(\"19\u00f78=2, 3\", \"38\u00f77=5, 3\"),\n    @(\"34\u00f79=3, 7\", \"71\u00f77=10, 1\"),\n    @(\"46\u00f74=11, 2\", \"83\u00f76=13, 5\"),\n    @(\"68\u00f74=17, 0\", \"16\u00f78=2, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
